# Fruta / hortaliza, semanal
# Insert two new rows of the latest weekly price-report data for
# "Zanahoria" at Terminal Hortofrutícola Agro Chillán (region Ñuble),
# pushing the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 391-392),
# shifting every existing row from 391 downward by two rows.
$ws.Rows("391:392").Insert()

# New row 391: "Primera" quality entry for the new reporting date.
$ws.Cells.Item(391, 1).Value  = 7
$ws.Cells.Item(391, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(391, 3).Value  = "Ñuble"
$ws.Cells.Item(391, 4).Value  = 45077
$ws.Cells.Item(391, 5).Value  = 16
$ws.Cells.Item(391, 6).Value  = 100114013
$ws.Cells.Item(391, 7).Value  = "Zanahoria"
$ws.Cells.Item(391, 8).Value  = "Sin especificar"
$ws.Cells.Item(391, 9).Value  = "Primera"
$ws.Cells.Item(391, 10).Value = 200
$ws.Cells.Item(391, 11).Value = 7000
$ws.Cells.Item(391, 12).Value = 7000
$ws.Cells.Item(391, 13).Value = 7000
$ws.Cells.Item(391, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(391, 15).Value = "Región de Ñuble"
$ws.Cells.Item(391, 16).Value = 350
$ws.Cells.Item(391, 17).Value = 20
$ws.Cells.Item(391, 18).Value = "Hortaliza"

# New row 392: "Segunda" quality entry for the same reporting date.
$ws.Cells.Item(392, 1).Value  = 7
$ws.Cells.Item(392, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(392, 3).Value  = "Ñuble"
$ws.Cells.Item(392, 4).Value  = 45077
$ws.Cells.Item(392, 5).Value  = 16
$ws.Cells.Item(392, 6).Value  = 100114013
$ws.Cells.Item(392, 7).Value  = "Zanahoria"
$ws.Cells.Item(392, 8).Value  = "Sin especificar"
$ws.Cells.Item(392, 9).Value  = "Segunda"
$ws.Cells.Item(392, 10).Value = 150
$ws.Cells.Item(392, 11).Value = 6000
$ws.Cells.Item(392, 12).Value = 6000
$ws.Cells.Item(392, 13).Value = 6000
$ws.Cells.Item(392, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(392, 15).Value = "Región de Ñuble"
$ws.Cells.Item(392, 16).Value = 300
$ws.Cells.Item(392, 17).Value = 20
$ws.Cells.Item(392, 18).Value = "Hortaliza"
